$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the data range B2:D9 to 0 (area check reset)
$ws.Range("B2:D9").Value = 0
